# Updates the cryptocurrency Price (D) and Volume/1h change (E) columns
# to match the latest scraped values from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.380.21'
$ws.Range('E2').Value = '  +1.34%  '
$ws.Range('D3').Value = '3.741.98'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'594.25"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').Value = "'166.82"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').Value = '3.740.15'
$ws.Range('E7').Value = '  -0.27%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -0.72%  '
$ws.Range('E10').Value = '  -2.52%  '
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').Value = "'0.0000258"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -5.41%  '
$ws.Range('D14').Value = "'36.08"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').Value = '4.369.60'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').Value = '3.788.80'
$ws.Range('E16').Value = '  +0.79%  '
$ws.Range('D17').Value = '68.388.60'
$ws.Range('E17').Value = '  +1.45%  '
$ws.Range('D18').Value = "'17.83"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.80%  '
$ws.Range('D19').Value = "'7.00"
$ws.Range('D19').ClearFormats()
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').Value = "'10.69"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.85%  '
$ws.Range('D22').Value = "'465.87"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.13%  '
$ws.Range('D23').Value = "'0.698"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.38%  '
$ws.Range('D24').Value = "'83.93"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').Value = "'0.0000144"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.95%  '
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').Value = "'12.03"
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.69%  '
$ws.Range('D28').Value = "'10.11"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.02%  '
$ws.Range('D30').Value = '3.888.62'
$ws.Range('E30').Value = '  -0.33%  '
$ws.Range('E31').Value = '  -4.26%  '
$ws.Range('D32').Value = "'7.28"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.28%  '
$ws.Range('D33').Value = "'29.83"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.64%  '
$ws.Range('D34').Value = "'2.18"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -1.86%  '
$ws.Range('D35').Value = "'9.21"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +1.38%  '
$ws.Range('D37').Value = '3.697.84'
$ws.Range('E37').Value = '  -0.42%  '
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('D39').Value = "'3.38"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -10.66%  '
$ws.Range('E40').Value = '  +0.23%  '
$ws.Range('D41').Value = "'0.997"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').Value = "'0.999"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('E45').Value = '  -1.92%  '
$ws.Range('E46').Value = '  -0.95%  '
$ws.Range('D47').Value = "'1.92"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('D48').Value = "'42.63"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +9.85%  '
$ws.Range('D49').Value = "'45.82"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').Value = "'146.25"
$ws.Range('D50').ClearFormats()
$ws.Range('D51').Value = "'389.28"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.55%  '
